$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 3x3 sending/target-cluster grid (ECs, FAPs, sCs) replaces the old
# 3-row FAPs/sCs-only table. Data below follows the Dr Hou-advised natmi
# recomputation.
$data = @(
    @("ECs","Sema3c","Nrp2","ECs", 2, 0.6666666666666666, 0.227933, 0.683799, 0.005274728560320285, 0.005274728560320286, 3, 1, 26.83081766666667, 80.49245300000001, 0.5916656861001716, 0.5916656861001716, 6.115628763216336, 55.04065886894701, 0.003120875892634072, 0.003120875892634072),
    @("ECs","Sema3c","Nrp2","FAPs", 2, 0.6666666666666666, 0.227933, 0.683799, 0.005274728560320285, 0.005274728560320286, 3, 1, 10.21969166666667, 30.659075, 0.2253617819930474, 0.2253617819930474, 2.329404980658334, 20.964644825925, 0.001188722227883401, 0.001188722227883401),
    @("ECs","Sema3c","Nrp2","sCs", 2, 0.6666666666666666, 0.227933, 0.683799, 0.005274728560320285, 0.005274728560320286, 3, 1, 8.297426666666667, 24.89228, 0.1829725319067811, 0.1829725319067811, 1.891257352413334, 17.02131617172, 0.0009651304398028127, 0.0009651304398028129),
    @("FAPs","Sema3c","Nrp2","ECs", 3, 1, 42.071953, 126.215859, 0.9736112457500787, 0.9736112457500787, 3, 1, 26.83081766666667, 80.49245300000001, 0.5916656861001716, 0.5916656861001716, 1128.82489982357, 10159.42409841213, 0.576052365711563, 0.576052365711563),
    @("FAPs","Sema3c","Nrp2","FAPs", 3, 1, 42.071953, 126.215859, 0.9736112457500787, 0.9736112457500787, 3, 1, 10.21969166666667, 30.659075, 0.2253617819930474, 0.2253617819930474, 429.9623874744917, 3869.661487270425, 0.2194147653107086, 0.2194147653107086),
    @("FAPs","Sema3c","Nrp2","sCs", 3, 1, 42.071953, 126.215859, 0.9736112457500787, 0.9736112457500787, 3, 1, 8.297426666666667, 24.89228, 0.1829725319067811, 0.1829725319067811, 349.0889447409467, 3141.80050266852, 0.1781441147278071, 0.1781441147278071),
    @("sCs","Sema3c","Nrp2","ECs", 3, 1, 0.912385, 2.737155, 0.02111402568960099, 0.021114025689601, 3, 1, 26.83081766666667, 80.49245300000001, 0.5916656861001716, 0.5916656861001716, 24.48003557680167, 220.320320191215, 0.01249244449597442, 0.01249244449597442),
    @("sCs","Sema3c","Nrp2","FAPs", 3, 1, 0.912385, 2.737155, 0.02111402568960099, 0.021114025689601, 3, 1, 10.21969166666667, 30.659075, 0.2253617819930474, 0.2253617819930474, 9.324293381291668, 83.91864043162501, 0.004758294454455461, 0.004758294454455462),
    @("sCs","Sema3c","Nrp2","sCs", 3, 1, 0.912385, 2.737155, 0.02111402568960099, 0.021114025689601, 3, 1, 8.297426666666667, 24.89228, 0.1829725319067811, 0.1829725319067811, 7.570447629266667, 68.1340286634, 0.003863286739171112, 0.003863286739171113)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $row[$j]
    }
}
